# Completed jobs; initialized the rest of preciseTAD jobs for gm12878 and k562
# using arrowhead for other resolutions (10kb-100kb); all total 2*4*21 jobs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("preciseTAD")

# --- 25kb block (rows 3-19, column H = k562_arrowhead @ 25kb) ---
# Mark chromosomes 1-8 (rows 3-10) as newly queued ("r") for the 25kb
# k562_arrowhead column; chromosome 2 (row 3/E3) and chromosome 5 (row
# 5/Q5) complete ("x" in GM12878 peakachu / k562 peakachu columns);
# chromosomes 11-19 (rows 11-19) that were queued ("r") are now done ("x").
$ws.Range("E3").Value = "x"
$ws.Range("H3").Value = "r"
$ws.Range("H4").Value = "r"
$ws.Range("H5").Value = "r"
$ws.Range("Q5").Value = "x"
for ($r = 6; $r -le 19; $r++) {
    $ws.Cells.Item($r, 8).Value = "x"
}

# --- Initialize the rest of preciseTAD jobs for gm12878 (col B) and
# k562 (col E) arrowhead, for resolutions 10kb, 25kb, 50kb, 100kb
# (row blocks 26-46, 49-69, 72-92, 95-115), marking each as running ("r").
$blockStarts = @(26, 49, 72, 95)
foreach ($start in $blockStarts) {
    $end = $start + 20
    for ($r = $start; $r -le $end; $r++) {
        $ws.Cells.Item($r, 2).Value = "r"
        $ws.Cells.Item($r, 5).Value = "r"
    }
}

# --- Update the saved selection to reflect the newly-filled range ---
$ws.Range("E95:E115").Select()
